$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: 10/03/2018 - Visionnage et lecture de tutos sur JavaFX - 3h
$ws.Range("A11").Value = 43169
$ws.Range("B11").Value = "Visionnage et lecture de tutos sur JavaFX"
$ws.Range("C11").Value = 3

# Row 12: 16/03/2018 - Discussion planification et organisation de groupe - 1.5h
$ws.Range("A12").Value = 43175
$ws.Range("B12").Value = "Discussion planification et organisation de groupe"
$ws.Range("C12").Value = 1.5

# Row 13: 17/03/2018 - Visionnage et lecture de tutos sur JavaFX et début de la création de l'interface graphique de la toolBar. - 4h
$ws.Range("A13").Value = 43176
$ws.Range("B13").Value = "Visionnage et lecture de tutos sur JavaFX et début de la création de l'interface graphique de la toolBar."
$ws.Range("C13").Value = 4
$ws.Rows("13:13").RowHeight = 30

# Row 14: 18/03/2018 - Visionnage et lecture de tutos sur JavaFX et suite de la création de l'interface graphique de la toolBar. - 4h
$ws.Range("A14").Value = 43177
$ws.Range("B14").Value = "Visionnage et lecture de tutos sur JavaFX et suite de la création de l'interface graphique de la toolBar."
$ws.Range("C14").Value = 4
$ws.Rows("14:14").RowHeight = 30

# Update the current selection to mirror the author's final cursor position
$ws.Range("G9").Select()

$wb.Save()
